$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Rename the categorical header labels (lower-cased) on every sheet /
#    table. Writing straight to the header cell also re-syncs the ListObject
#    (table) column name automatically.
# ---------------------------------------------------------------------------
$headerMap = @{
    "ID"           = "id_global"
    "Departamento" = "departamento"
    "Municipio"    = "municipio"
    "Colegio"      = "colegio"
    "Sede"         = "sede"
    "Género"       = "genero"
    "Grado"        = "grado"
    "Grupo"        = "grupo"
    "Jornada"      = "jornada"
}

foreach ($ws in $wb.Worksheets) {
    $lastCol = $ws.UsedRange.Columns.Count
    for ($c = 1; $c -le $lastCol; $c++) {
        $cell = $ws.Cells.Item(1, $c)
        $cur = $cell.Value()
        if ($headerMap.ContainsKey($cur)) {
            $cell.Value = $headerMap[$cur]
        }
    }
}

# ---------------------------------------------------------------------------
# 2. Select A1:L1 (the renamed categorical columns) on every sheet.
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)
$ws3 = $wb.Worksheets.Item(3)
$ws4 = $wb.Worksheets.Item(4)

$ws1.Range("A1:L1").Select()
$ws2.Range("A1:L1").Select()
$ws3.Range("A1:L1").Select()
$ws4.Range("A1:L1").Select()

# ---------------------------------------------------------------------------
# 3. New "Highlight duplicates" rule on A1 (blue fill) for sheets 2-4, and
#    shrink the pre-existing whole-column rule so it starts at row 2 where
#    it still covers the whole column (sheet 1 keeps its original,
#    untouched rule).
# ---------------------------------------------------------------------------
function Add-HeaderDuplicateRule($ws) {
    $rngA1 = $ws.Range("A1")
    $fc = $rngA1.FormatConditions.AddUniqueValues(1)
    $fc.DupeUnique = 1
    $fc.Interior.Color = 12611584   # OLE BGR for FF0070C0
    $fc.SetFirstPriority()
}

# Sheet 2 ("3º"): shrink the full-column rule to A2:A1048576, then add A1 rule
$fcs2 = $ws2.Range("A1:A1048576").FormatConditions
$fcs2.Item(1).ModifyAppliesToRange($ws2.Range("A2:A1048576"))
Add-HeaderDuplicateRule $ws2

# Sheet 3 ("4º"): only add the new A1 rule (existing A2 rule untouched)
Add-HeaderDuplicateRule $ws3

# Sheet 4 ("5º"): shrink the full-column rule to A2:A1048576, then add A1 rule
$fcs4 = $ws4.Range("A1:A1048576").FormatConditions
$fcs4.Item(1).ModifyAppliesToRange($ws4.Range("A2:A1048576"))
Add-HeaderDuplicateRule $ws4

# ---------------------------------------------------------------------------
# 4. Active sheet becomes "5º" (4th tab, 0-based activeTab = 3).
# ---------------------------------------------------------------------------
$ws4.Activate()
